$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Chart titles: "Cycle GAN" -> "CycleGAN" rebrand (two charts on the sheet)
# ---------------------------------------------------------------------------
$co1 = $ws.ChartObjects().Item(1)
$chart1 = $co1.Chart
$title1 = $chart1.ChartTitle
$title1.Text = "Loss of Discriminator_A for different types of CycleGAN"

$co2 = $ws.ChartObjects().Item(2)
$chart2 = $co2.Chart
$title2 = $chart2.ChartTitle
$title2.Text = "Loss of Generator_A_to_B for different types of CycleGAN"

# Give chart2's title a manual (edge-anchored) layout position.
$title2.Left = 0.16308314550214639
$title2.Top = 0

# ---------------------------------------------------------------------------
# 2) Resize / reposition the two chart frames anchored on the worksheet
# ---------------------------------------------------------------------------
# Chart 1 ("Diagramm 2"): top-left anchor is unchanged, only grows a bit.
$co1.Width = 465.4
$co1.Height = 328.42503937007876

# Chart 2 ("Diagramm 3"): shifts up a hair and grows taller/narrower.
$co2.Top = 198.4499212598425
$co2.Width = 461.0500787401575
$co2.Height = 329.54999999999995

# ---------------------------------------------------------------------------
# 3) Update the sheet view: drop the frozen top-left scroll position and move
#    the active selection from I18 to I40.
# ---------------------------------------------------------------------------
$ws.Range("I40").Select() | Out-Null
